$d = $word.ActiveDocument

function Replace-All($old, $new) {
    $keepGoing = $true
    while ($keepGoing) {
        $rng = $d.Content
        $found = $rng.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $null, 0)
        if ($found) {
            $rng.Text = $new
        } else {
            $keepGoing = $false
        }
    }
}

Replace-All "Important message to deliver: the composition of substances can affect the appearance and properties of objects. The change in composition can manifest itself as a change in the object appearance" "Importante messaggio da trasmettere: la composizione delle sostanze può influenzare l'aspetto e le proprietà degli oggetti. Il cambiamento della composizione può manifestarsi nell'aspetto dell'oggetto"
Replace-All "Note: pepper or coffee are 'spectators' of the transformation, they serve only to visualize the change in surface tension. The transformation depends only on the soap addition to water." "Nota: il pepe o il caffè sono 'spettatori' della trasformazione, servono solo per vedere il cambiamento nella tensione superficiale. La trasformazione dipende solo dall'aggiunta di sapone all'acqua."
Replace-All "Set for each group of students: a glass or a cup, a plate, one stick, soap, water, coffee, black pepper. Both the plates and the water should be clean." "Set per ogni gruppo di studenti: un bicchiere o una tazza, un piatto, un bastone, sapone, acqua, caffè, pepe nero. Sia l'acqua che i piatti dovrebbero esser puliti."
Replace-All "Because of their chemical nature, some substances tend to concentrate in specific regions, while some others simply spread randomly" "A causa della loro natura chimica, alcune sostanze tendono a concentrarsi in regioni specifiche, mentre alcune altre, semplicemente, si diffondono casualmente"
Replace-All "Facilitate the discussion: why are the coffee/pepper particles pushed to the plate edge? What is the role of soap?" "Facilitare la discussione: perché le particelle di caffè/pepe sono spinte al bordo del piatto? Qual è il ruolo del sapone?"
Replace-All "Discuss with the other learners why does this effect occur and if they have noticed it before in  daily life." "Discutono con gli altri studenti perché si verifica quest'effetto e se lo hanno notato prima nella loro vita quotidiana."
Replace-All "As coffees are mixtures and their composition can vary, some coffee can react poorly during the experiment." "Poiché i caffè sono miscugli e la loro composizione può variare, alcuni caffè potrebbero reagire male durante l'esperimento."
Replace-All "Pour some black pepper or coffee on top of the water. Observe the uniform distribution of the particles" "Versano del pepe nero o caffè sull'acqua. Osservano la distribuzione uniforme delle particelle"
Replace-All "After the experiment, the distribution of particles cannot be further modified with the same method." "Dopo l'esperimento, la distribuzione delle particelle non è ulteriormente modificabile con lo stesso metodo."
Replace-All "Observe the spreading of the particles, or, generally, the modification of their distribution" "Osservano la diffusione delle particelle, o, in generale, la modifica della loro distribuzione"
Replace-All "Learn the nature of surface tension in water and its modifications with everyday objects." "Scoprire la natura della tensione superficiale in acqua e i suoi cambiamenti con oggetti d'uso quotidiano."
Replace-All "Suggestion for discussion: surface tension depends on water surface composition " "Suggerimenti per la discussione: la tensione superficiale dipende dalla composizione della superficie dell'acqua "
Replace-All "Suggestion for discussion: surface tension is a surface property" "Suggerimento per la discussione: la tensione superficiale è una proprietà della superficie"
Replace-All "Observe the effects of surface tension (curved water surface). " "Osservano l'effetto della tensione superficiale (superficie curva dell'acqua). "
Replace-All "To repeat the experiment, first, clean the plate thoroughly." "Per ripetere l'esperimento, pulire prima il piatto in modo accurato."
Replace-All "Fill the glasses or cups up to the very top" "Riempiono i bicchieri o le tazze fino in cima"
Replace-All "Put a drop of soap on the tip of the sticks" "Versano una goccia di sapone sulla punta dei bastoncini"
Replace-All "Put a layer of water on top of the plates" "Mettono uno strato d'acqua sui piatti"
Replace-All "Touch the water surface with the stick" "Toccano la superfici d'acqua con il bastone"
Replace-All "Introduction of the second experiment" "Introduzione al secondo esperimento"
Replace-All "Introduction of the first experiment" "Introduzione al primo esperimento"
Replace-All "Assist the process, provoke thoughts" "Assiste il processo, causa pensieri"
Replace-All "Try out guesses and share ideas " "Provano a indovinare e condividono le proprie idee "
Replace-All "General VMC Video Introduction" "Introduzione Generale al Video di VMC"
Replace-All "Filling the glasses with water" "Riempire i bicchieri d'acqua"
Replace-All "Why do the particles spread?" "Perché le particelle si diffondono?"
Replace-All "Experiment solution (part 1)" "Soluzione all'esperimento (parte 1)"
Replace-All "Experiment solution (part 2)" "Soluzione all'esperimento (parte 2)"
Replace-All "Modifying surface tension" "Modificare la tensione superficiale"
Replace-All "Facilitate the discussion" "Facilita la discussione"
Replace-All "Invitation to discussion" "Invito alla discussione"
Replace-All "Where does the soap go?" "Dove finisce il sapone?"
Replace-All "What facilitator does" "Cosa fa il facilitatore"
Replace-All "Video Introduction" "Introduzione al Video"
Replace-All "What learners do" "Cosa fanno gli studenti"
Replace-All "Surface tension" "Tensione superficiale"
Replace-All "N. of students" "N. di studenti"
Replace-All "Camp Location" "Posizione del Campo"
Replace-All "Facilitators" "Facilitatori"
Replace-All "Preparations" "Preparazioni"
Replace-All "Video Title" "Titolo del Video"
Replace-All "VIDEO PAUSE" "VIDEO IN PAUSA"
Replace-All "Experiment:" "Esperimento:"
Replace-All "Discussion:" "Discussione:"
Replace-All "Video time" "Tempo del video"
Replace-All "Conclusion" "Conclusione"
Replace-All "Chemistry" "Chimica"
Replace-All "Resources" "Risorse"
Replace-All "Comments" "Commenti"
Replace-All "Material" "Materiale"
Replace-All "Aim(s)" "Obiettivo/i"
Replace-All "Length" "Lunghezza"
Replace-All "needed" "necessarie"
Replace-All "Topic" "Argomento"
Replace-All "Date" "Data"
Replace-All "None" "Nessuna"
